$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: was M, now B
$ws.Range("A2").Value = "B"
$ws.Range("B2").Value = 0.9102564102564102
$ws.Range("C2").Value = 0.9861111111111112
$ws.Range("D2").Value = 0.9466666666666667
$ws.Range("E2").Value = 72

# Row 3: was B, now M
$ws.Range("A3").Value = "M"
$ws.Range("B3").Value = 0.9722222222222222
$ws.Range("C3").Value = 0.8333333333333334
$ws.Range("D3").Value = 0.8974358974358975
$ws.Range("E3").Value = 42

# Row 4: accuracy
$ws.Range("B4").Value = 0.9298245614035088
$ws.Range("C4").Value = 0.9298245614035088
$ws.Range("D4").Value = 0.9298245614035088
$ws.Range("E4").Value = 0.9298245614035088

# Row 5: macro avg
$ws.Range("B5").Value = 0.9412393162393162
$ws.Range("C5").Value = 0.9097222222222223
$ws.Range("D5").Value = 0.9220512820512821

# Row 6: weighted avg
$ws.Range("B6").Value = 0.9330859199280251
$ws.Range("C6").Value = 0.9298245614035088
$ws.Range("D6").Value = 0.9285290148448043
